$p = $ppt.ActivePresentation

# Identify the title text of a slide (first shape whose name starts with "Title").
function Get-SlideTitleText($slide) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Name -like "Title*" -and $sh.HasTextFrame) {
            return $sh.TextFrame.TextRange.Text
        }
    }
    return $null
}

# Remove the slides that should no longer be part of the deck: "Capstone Project"
# (old title slide), "MVP target", "Technical Achievement" and "Questions?".
# Walk from the last slide to the first so deleting doesn't disturb the
# indices of slides still to be inspected.
$titlesToRemove = @("Capstone Project", "MVP target", "Technical Achievement", "Questions?")
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $title = Get-SlideTitleText $slide
    if ($titlesToRemove -contains $title) {
        $slide.Delete()
    }
}

# Remaining slides, in order, are now:
#   1: Background (brotherhood...)
#   2: Background: part 2
#   3: Profile Page / Home Page / ... mock-ups, previously titled "End Of The Year Product Target"
#   4: Reference

# On the mock-ups slide, reposition/resize the title textbox and retitle it.
$s = $p.Slides.Item(3)
$title = $s.Shapes.Item("Title 1")
$title.Left = 4303471 / 914400 * 72
$title.Top = 88663 / 914400 * 72
$title.Width = 5178580 / 914400 * 72
$title.Height = 1325563 / 914400 * 72
$title.TextFrame.TextRange.Text = "Project description"
